$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 273840.5
$ws.Range("C4").Value = 0.6543422341346741
